# Applies the "added sample json file" commit's spreadsheet edits:
#  - A1 header text gets a clarifying suffix
#  - E1 header renamed
#  - a new "KubeSecurityEvents" sample row is appended at row 20
#  - column A is widened to fit the new, longer header text
#  - the active selection moves to F1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text tweaks (row 1) ---
$ws.Range("A1").Value = "ClusterLevel(clusterName)"
$ws.Range("E1").Value = "policyLastUpdated"

# --- New sample block: KubeSecurityEvents ---
$ws.Range("A20").Value = "KubeSecurityEvents"
$ws.Range("B20").Value = "EventName"
$ws.Range("C20").Value = "LoggedonUser"
$ws.Range("D20").Value = "TimeGen"

# --- Column A autosize to fit the new, longer header ---
$ws.Columns.Item(1).ColumnWidth = 21.833333333333332

# --- Move the active selection to F1 ---
$ws.Range("F1").Select()
